$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.973.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.080.67'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.078.81'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.451'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.396'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.608.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000161'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.033.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.075.63'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '347.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.73%  '
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.498'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.32%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0863'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.39%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.87'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.93'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.05%  '
$ws.Range('E36').Value = '  -3.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.96'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  -1.59%  '
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.58'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.689'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.385.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '36.59'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.80%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.118.56'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('E50').Value = '  -3.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.18%  '
